{"js": "// Replace the 100 arithmetic-problem cell values in the single 20x5 table\n// with the values from the updated worksheet (row-major order, matching\n// the order the cells appear in the document).\nconst newValues = [\n  [\"59+7=\", \"47+7=\", \"19+64=\", \"77+4=\", \"75-36=\"],\n  [\"17+79=\", \"57+36=\", \"48-9=\", \"5+67=\", \"59+3=\"],\n  [\"81-75=\", \"57+26=\", \"78+13=\", \"20-15=\", \"2+59=\"],\n  [\"66-59=\", \"65-27=\", \"83-78=\", \"53-24=\", \"81-18=\"],\n  [\"6+86=\", \"63-25=\", \"50-9=\", \"14+78=\", \"23+59=\"],\n  [\"42-13=\", \"28+45=\", \"18+28=\", \"77+17=\", \"73+18=\"],\n  [\"64+8=\", \"40-29=\", \"62-58=\", \"6+78=\", \"54-38=\"],\n  [\"50-19=\", \"93-88=\", \"97-39=\", \"33+58=\", \"98-89=\"],\n  [\"76+9=\", \"17+74=\", \"44-6=\", \"71-28=\", \"98-39=\"],\n  [\"41-8=\", \"61-38=\", \"80-28=\", \"47+24=\", \"7+79=\"],\n  [\"50-27=\", \"72-63=\", \"77+19=\", \"59+35=\", \"42-8=\"],\n  [\"9+23=\", \"91-55=\", \"74-38=\", \"66-58=\", \"39+45=\"],\n  [\"97-38=\", \"47+19=\", \"8+37=\", \"40-22=\", \"27-18=\"],\n  [\"8+78=\", \"31-17=\", \"25+39=\", \"64+8=\", \"9+83=\"],\n  [\"64-39=\", \"18+47=\", \"59+32=\", \"28+37=\", \"45-38=\"],\n  [\"63+28=\", \"19+49=\", \"17+74=\", \"8+55=\", \"22-18=\"],\n  [\"76+16=\", \"8+8=\", \"80-39=\", \"90-47=\", \"57+17=\"],\n  [\"50-49=\", \"51-18=\", \"7+77=\", \"9+53=\", \"35+37=\"],\n  [\"72-24=\", \"73-17=\", \"18+79=\", \"63+8=\", \"71-48=\"],\n  [\"32+9=\", \"64-28=\", \"19+52=\", \"77+5=\", \"19+53=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r][c];\n  }\n}\nawait context.sync();\n\n", "ps1": "# Replace the 100 arithmetic-problem cell values in the single 20x5 table\n# with the values from the updated worksheet (row-major order, matching\n# the order the cells appear in the document).\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"59+7=\", \"47+7=\", \"19+64=\", \"77+4=\", \"75-36=\"),\n    @(\"17+79=\", \"57+36=\", \"48-9=\", \"5+67=\", \"59+3=\"),\n    @(\"81-75=\", \"57+26=\", \"78+13=\", \"20-15=\", \"2+59=\"),\n    @(\"66-59=\", \"65-27=\", \"83-78=\", \"53-24=\", \"81-18=\"),\n    @(\"6+86=\", \"63-25=\", \"50-9=\", \"14+78=\", \"23+59=\"),\n    @(\"42-13=\", \"28+45=\", \"18+28=\", \"77+17=\", \"73+18=\"),\n    @(\"64+8=\", \"40-29=\", \"62-58=\", \"6+78=\", \"54-38=\"),\n    @(\"50-19=\", \"93-88=\", \"97-39=\", \"33+58=\", \"98-89=\"),\n    @(\"76+9=\", \"17+74=\", \"44-6=\", \"71-28=\", \"98-39=\"),\n    @(\"41-8=\", \"61-38=\", \"80-28=\", \"47+24=\", \"7+79=\"),\n    @(\"50-27=\", \"72-63=\", \"77+19=\", \"59+35=\", \"42-8=\"),\n    @(\"9+23=\", \"91-55=\", \"74-38=\", \"66-58=\", \"39+45=\"),\n    @(\"97-38=\", \"47+19=\", \"8+37=\", \"40-22=\", \"27-18=\"),\n    @(\"8+78=\", \"31-17=\", \"25+39=\", \"64+8=\", \"9+83=\"),\n    @(\"64-39=\", \"18+47=\", \"59+32=\", \"28+37=\", \"45-38=\"),\n    @(\"63+28=\", \"19+49=\", \"17+74=\", \"8+55=\", \"22-18=\"),\n    @(\"76+16=\", \"8+8=\", \"80-39=\", \"90-47=\", \"57+17=\"),\n    @(\"50-49=\", \"51-18=\", \"7+77=\", \"9+53=\", \"35+37=\"),\n    @(\"72-24=\", \"73-17=\", \"18+79=\", \"63+8=\", \"71-48=\"),\n    @(\"32+9=\", \"64-28=\", \"19+52=\", \"77+5=\", \"19+53=\"),\n)\n\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Length; $c++) {\n        $cell = $tbl.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $row[$c]\n    }\n}\n"}
